$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.765.06'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.32%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.507.93'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.47%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '608.94'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '192.16'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.35%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.629'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.03%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.212'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.93%  '
$ws.Range('E10').Value = '  +2.68%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.50'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.26%  '
$ws.Range('E12').Value = '  -2.81%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.60'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.22%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.066.84'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.43%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '617.53'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +7.86%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '69.860.64'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.97'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.56%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.63'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.47%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.518.37'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.94%  '
$ws.Range('E20').Value = '  -0.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.990'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.07%  '
$ws.Range('B22').Value = 'InternetComputer(DFINITY)'
$ws.Range('C22').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '17.53'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.61%  '
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '107.28'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +14.00%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.65'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.56%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.04'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.81%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.05'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.29%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.98'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.32%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.77'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.53%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '34.25'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +5.54%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.98'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.04%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '12.43'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.71%  '
$ws.Range('B32').Value = 'dogwifhat'
$ws.Range('C32').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.06'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +6.47%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.115'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.83'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.56%  '
$ws.Range('B35').Value = 'Maker'
$ws.Range('C35').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.691.62'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.31%  '
$ws.Range('B36').Value = 'Fetch.AI'
$ws.Range('C36').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.07'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.82%  '
$ws.Range('B37').Value = 'Dai'
$ws.Range('C37').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('B38').Value = 'Bittensor'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '514.40'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.55%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.60'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.91%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.391'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.17%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.72'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.37%  '
$ws.Range('E42').Value = '  -2.56%  '
$ws.Range('E43').Value = '  -1.25%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0465'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.51%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.89'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.46%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.142'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.58%  '
$ws.Range('E47').Value = '  -4.25%  '
$ws.Range('B48').Value = 'FirstDigitalUSD'
$ws.Range('C48').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.00'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.47%  '
$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.71'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -5.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '132.51'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.29%  '
$ws.Range('B51').Value = 'OceanProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.34'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -6.25%  '
